$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (apostrophe-prefixed so Excel
# stores it as literal text, matching the original t="inlineStr" cells,
# instead of auto-converting numeric-/percent-looking text into a real
# number or percentage).
$changes = [ordered]@{
  "D2" = "'307.19"
  "E2" = "'1.84%"
  "D3" = "'35.82"
  "E3" = "'0.77%"
  "D4" = "'5.039"
  "E4" = "'-0.16%"
  "E5" = "'1.57%"
  "D6" = "'1.956"
  "E6" = "'1.77%"
  "D7" = "'4.131"
  "E7" = "'1.65%"
  "D8" = "'7.795"
  "E8" = "'-0.20%"
  "D9" = "'0.9325"
  "E9" = "'0.80%"
  "D10" = "'0.1341"
  "E10" = "'-6.58%"
  "D11" = "'0.1925"
  "E11" = "'1.45%"
  "D12" = "'0.09247"
  "E12" = "'-0.01%"
  "D13" = "'0.03508"
  "E13" = "'3.45%"
  "D14" = "'0.09875"
  "E14" = "'0.05%"
  "D15" = "'0.001434"
  "E15" = "'3.05%"
  "D16" = "'0.005826"
  "E16" = "'1.07%"
  "D17" = "'3.602"
  "E17" = "'2.41%"
  "D18" = "'2.942"
  "E18" = "'-0.57%"
  "D19" = "'0.3428"
  "E19" = "'0.71%"
  "D20" = "'0.1332"
  "E20" = "'2.21%"
  "D21" = "'5.211"
  "E21" = "'3.08%"
  "E22" = "'7.99%"
  "D23" = "'0.04387"
  "E23" = "'-2.30%"
  "E24" = "'0.32%"
  "D25" = "'0.004775"
  "E25" = "'-0.29%"
  "E26" = "'34.13%"
  "E27" = "'3.98%"
  "D39" = "'0.01988"
  "E39" = "'3.83%"
  "D40" = "'0.05009"
  "E40" = "'5.51%"
  "D41" = "'0.01120"
  "E41" = "'15.84%"
  "D42" = "'0.007627"
  "E42" = "'3.75%"
  "D43" = "'0.1379"
  "E43" = "'3.87%"
  "D44" = "'0.002096"
  "E44" = "'-0.69%"
  "D45" = "'0.01133"
  "E45" = "'9.86%"
  "D46" = "'0.00006385"
  "E46" = "'2.07%"
  "E47" = "'-0.26%"
  "E49" = "'-28.44%"
  "E50" = "'-0.26%"
  "E51" = "'-0.26%"
}

foreach ($addr in $changes.Keys) {
  $cell = $ws.Range($addr)
  $cell.Value = $changes[$addr]
  # Reset to the Normal style so the quote-prefix text entry does not
  # leave a stray cell style behind (original cells used the default style).
  $cell.Style = "Normal"
}
